# C1: Add a "Driver_Settings" sheet to the example catdriver_config workbook
# and drop the leftover empty "Order" placeholder cells on the Variables sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheet at the end of the workbook (after "Variables")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Driver_Settings"

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "driver"
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "levels_order"
$ws.Range("D1").Value = "reference_level"
$ws.Range("E1").Value = "missing_strategy"
$ws.Range("F1").Value = "rare_level_policy"

# Bold, white text on a solid blue fill, centered horizontally
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 12874308
$headerRange.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------------
# grade: ordinal, D;C;B;A order, D as reference level
$ws.Range("A2").Value = "grade"
$ws.Range("B2").Value = "ordinal"
$ws.Range("C2").Value = "D;C;B;A"
$ws.Range("D2").Value = "D"
$ws.Range("E2").Value = "missing_as_level"
$ws.Range("F2").Value = "warn_only"

# campus: categorical
$ws.Range("A3").Value = "campus"
$ws.Range("B3").Value = "categorical"
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = "missing_as_level"
$ws.Range("F3").Value = "warn_only"

# course_type: categorical
$ws.Range("A4").Value = "course_type"
$ws.Range("B4").Value = "categorical"
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = "missing_as_level"
$ws.Range("F4").Value = "warn_only"

# employment_field: categorical
$ws.Range("A5").Value = "employment_field"
$ws.Range("B5").Value = "categorical"
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = "missing_as_level"
$ws.Range("F5").Value = "warn_only"

# ---------------------------------------------------------------------------
# 4. Column widths (offset compensates for Excel's automatic padding so the
#    stored <col width> ends up exactly 19,14,15,18,19,20)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws.Columns.Item(2).ColumnWidth = 13.166666666666668
$ws.Columns.Item(3).ColumnWidth = 14.166666666666668
$ws.Columns.Item(4).ColumnWidth = 17.166666666666668
$ws.Columns.Item(5).ColumnWidth = 18.166666666666668
$ws.Columns.Item(6).ColumnWidth = 19.166666666666668

# ---------------------------------------------------------------------------
# 5. Clean up the now-unused empty "Order" inline-string placeholder cells
#    on the Variables sheet (D4, D5, D6 — Campus/Course Type/Employment
#    Field rows never had an Order value).
# ---------------------------------------------------------------------------
$variables = $wb.Worksheets.Item("Variables")
$variables.Range("D4").ClearContents()
$variables.Range("D5").ClearContents()
$variables.Range("D6").ClearContents()

# ---------------------------------------------------------------------------
# 6. Restore the original active sheet/selection (Settings!A1) so the only
#    changes are the new sheet and the Variables cleanup.
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$null = $settings.Activate()
$null = $settings.Range("A1").Select()

Write-Output "Driver_Settings sheet added; Variables sheet cleaned up."
